$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (A1/B1): "COD_DEPTO" -> "Valor", "Departamento" -> "Categoría"
$ws.Range("A1").Value = "Valor"
$ws.Range("B1").Value = "Categoría"

# Center-align the header row
$ws.Range("A1:B1").HorizontalAlignment = -4108

# Clear the stored cell selection (B26) so the sheet view has no selection override
$ws.Range("A1").Select()
